{"js": "// Insert the literal text \"<rename>\" (coloured blue, RGB 0070C0) immediately\n// after the first occurrence of \"HackDesk\" in the \"Sticky Boundary\" paragraph:\n//\n//   Before: Between the HackDesk and the Album screens ...\n//   After:  Between the HackDesk<rename> and the Album screens ...\n//\n// We search the body for \"HackDesk\" (case sensitive, whole document) and use\n// the occurrence that lives in the paragraph that also contains\n// \"Sticky Boundary\" / \"Between the\" (the first occurrence in the document),\n// so the later, unrelated \"HackDesk\" inside the \"Common Edit Mode\" paragraph\n// is left untouched.\n\nconst body = context.document.body;\nconst results = body.search(\"HackDesk\", { matchCase: true, matchWholeWord: true });\nresults.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < results.items.length; i++) {\n  const item = results.items[i];\n  const para = item.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n  if (para.text.indexOf(\"Between the\") !== -1) {\n    target = item;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find \"HackDesk\" in the \"Sticky Boundary\" paragraph.');\n}\n\n// Insert the new run right after the matched \"HackDesk\" text and colour it.\nconst inserted = target.insertText(\"<rename>\", Word.InsertLocation.after);\ninserted.font.color = \"#0070C0\";\n\nawait context.sync();\n", "ps1": "# Insert the literal text \"<rename>\" (coloured blue, RGB 0070C0) immediately\n# after the first occurrence of \"HackDesk\" in the \"Sticky Boundary\" paragraph:\n#\n#   Before: Between the HackDesk and the Album screens ...\n#   After:  Between the HackDesk<rename> and the Album screens ...\n#\n# There are two occurrences of \"HackDesk\" in the document (the other one is\n# in the \"Common Edit Mode\" paragraph, \"... from HackDesk or Album screens)\n# ...\"); we only want the first one, inside the paragraph that also contains\n# \"Sticky Boundary:\" / \"Between the\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"HackDesk\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop\n\n$target = $null\nwhile ($find.Execute()) {\n    $candidate = $find.Parent\n    $paraText = $candidate.Paragraphs(1).Range.Text\n    if ($paraText -like \"*Between the*\") {\n        $target = $d.Range($candidate.Start, $candidate.End)\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find 'HackDesk' inside the Sticky Boundary paragraph.\"\n}\n\n# Collapse to the end of the matched word and insert the new run right there.\n$target.Collapse(0)  # wdCollapseEnd\n$target.InsertAfter(\"<rename>\")\n$target.Font.Color = 12611584  # RGB(0x00, 0x70, 0xC0) -> 0070C0\n"}
